$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C (PartnerAlias) entirely, shifting remaining columns left
$ws.Range("C:C").Delete()

# Update active selection to match target state
$ws.Range("E6").Select()
